$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Create the new "KidsAlt" paragraph style (based on "Kids") before it
#    gets referenced anywhere, so Word doesn't silently auto-mint a
#    generic placeholder style for us.
# ---------------------------------------------------------------------
$kidsAlt = $d.Styles.Add("KidsAlt", 1)
$kidsAlt.BaseStyle = $d.Styles("Kids")
$kidsAlt.QuickStyle = $true
# Kids style indent is 936 twips (46.8 pt) left / hanging; KidsAlt moves
# that out to 1080 twips (54 pt) left / hanging, and adjusts its tab
# stops (right tab at 936twips/46.8pt, left tab at 1080twips/54pt),
# clearing the two inherited tab stops from "Kids" (576, 720 twips).
$kidsAlt.ParagraphFormat.LeftIndent = 54
$kidsAlt.ParagraphFormat.FirstLineIndent = -54
$kidsAlt.ParagraphFormat.TabStops.ClearAll()
$kidsAlt.ParagraphFormat.TabStops.Add(46.8, 2)
$kidsAlt.ParagraphFormat.TabStops.Add(54, 0)

# ---------------------------------------------------------------------
# 2) Re-style the "Fourth Child" paragraph from "Kids" to "KidsAlt".
# ---------------------------------------------------------------------
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "3`tv.`tFourth Child. Vital statistics here.") {
        $para.Style = $kidsAlt
        break
    }
}

# ---------------------------------------------------------------------
# 3) Remove the "Added for Ged2Reg." paragraph (styled "ReportInfo").
# ---------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $para = $d.Paragraphs($i)
    if ($para.Range.Text.TrimEnd([char]13, [char]7) -eq "Added for Ged2Reg.") {
        $para.Range.Delete()
        break
    }
}

# ---------------------------------------------------------------------
# 4) Drop the "ReportInfo" style definition (now unused).
# ---------------------------------------------------------------------
$d.Styles("ReportInfo").Delete()

# ---------------------------------------------------------------------
# 5) Merge the two footnote runs that straddle the lastRenderedPageBreak
#    into a single run (removing the stale rendered-page-break marker).
# ---------------------------------------------------------------------
$openQuote = [char]0x201C
$closeQuote = [char]0x201D
$footnoteText = "(Boston: NEHGS, 2014). Thanks to Helen Schatvet Ullmann, who inspired these instructions with her helpful online article " + $openQuote + "Register Style Template: A Template and Suggestions for Writing in Register Style in Microsoft Word," + $closeQuote + " online at AmericanAncestors.org/register-template. This footnote text style is applied automatically."

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute($footnoteText, $true, $false, $false, $false, $false, $true, 1, $false, $footnoteText, 2)

Write-Host "Edit complete"
